$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "I am writing to express my interest in the Senior Analyst - Data Science position at Tiger Analytics. With my strong background in data science and machine learning, I believe I would be a valuable addition to your team.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I am writing to express my interest in the Senior Associate Data Scientist position at your company. With a strong background in data science and a proven track record of delivering impactful insights and solutions, I believe I would be a valuable addition to your team.",
    2
) | Out-Null

$d.Content.Find.Execute(
    "I have over 4 years of experience in the field of data science, with a focus on developing and implementing advanced analytics solutions. In my previous role as a Senior Associate Data Scientist at Affine, I led a team in scraping data from multiple gaming sites and created concise dashboards with insightful charts using Google Sheets. I also supervised a data science team in generating and clustering a large dataset of players, and provided real-time visibility into retention and revenue predictions using FB-Prophet. Additionally, I have experience in analyzing sentiments from web servers and designing live revenue models for upcoming games.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In my current role as a Senior Associate Data Scientist at ffine, I have led a team in scraping data from multiple gaming sites and creating concise dashboards with insightful charts using Google Sheets. I have also supervised a data science team in generating and clustering analytical datasets, providing real-time visibility to leadership for decision-making. Additionally, I have analyzed sentiments from web servers and designed live revenue models and dashboards for game launches.",
    2
) | Out-Null

$d.Content.Find.Execute(
    "I have a comprehensive knowledge of SQL and am proficient in Python, with a passion for writing high-quality, modular, and scalable code. I am adept at using various data science approaches, machine learning algorithms, and statistical methods. I also have experience in customer analytics and marketing analytics, which I believe would be beneficial in this role.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "During my internship at iNeuron, I worked on an insurance premium prediction project where I conducted exploratory data analysis, implemented regression models, and deployed the project as a web application on AWS. I also have experience as an Assistant System Engineer at Tata Consultancy Services, where I performed system analysis, documentation, and user support.",
    2
) | Out-Null

$d.Content.Find.Execute(
    "I am a self-motivated and detail-oriented individual, with excellent problem-solving and analytical skills. I am able to effectively communicate and collaborate with stakeholders across the organization. I am also committed to continuous learning and staying updated with the latest advancements in the field of data science.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I hold a B.Tech degree from NIT Durgapur and I am currently pursuing a Full Stack Data Science Bootcamp at iNeuron. I have completed coursework in machine learning, probability and statistics, calculus, and stock market analysis. I am proficient in Python, SQL, NLP, time series forecasting, deep learning, and various data science libraries such as NumPy, Pandas, TensorFlow, and Keras.",
    2
) | Out-Null

$d.Content.Find.Execute(
    "I am excited about the opportunity to work at Tiger Analytics and contribute to solving complex problems using data and technology. I believe my skills and experience make me a strong fit for this role, and I am confident in my ability to make a positive impact on your team.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In addition to my technical skills, I possess strong soft skills including time management, teamwork, problem-solving, and presentation skills. I am a hardworking and flexible individual who is always eager to learn and adapt to new challenges.",
    2
) | Out-Null

$d.Content.Find.Execute(
    "Thank you for considering my application. I look forward to the opportunity to discuss how my skills and experience align with your needs.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I am excited about the opportunity to contribute my skills and experience to your team. Thank you for considering my application. I look forward to the possibility of discussing how I can contribute to your organization.",
    2
) | Out-Null

Write-Output "done"
